$wb = $excel.ActiveWorkbook

# tot-arrecad: fill in the missing UF label on the stray row 25
$wsTot = $wb.Worksheets.Item("tot-arrecad")
$wsTot.Range("A25").Value = "XX"

# avg-arrecad: same fix
$wsAvg = $wb.Worksheets.Item("avg-arrecad")
$wsAvg.Range("A25").Value = "XX"

# max-arrecad: same fix
$wsMax = $wb.Worksheets.Item("max-arrecad")
$wsMax.Range("A25").Value = "XX"

# tx-sucesso: fill/shift the UF labels for rows 3-5
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A3").Value = "MT"
$wsTx.Range("A4").Value = "MA"
$wsTx.Range("A5").Value = "XX"
